$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" values that must remain stored as text
# (matching the source data which keeps values like "311.55" or "1.001" as strings).
# Force text format first so Excel does not auto-convert them to numbers.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the numeric-looking text values
$ws.Range("D5").Value = '311.55'
$ws.Range("D7").Value = '0.5068'
$ws.Range("D8").Value = '0.3934'
$ws.Range("D9").Value = '0.09706'
$ws.Range("D10").Value = '1.142'
$ws.Range("D11").Value = '40.89'
$ws.Range("D12").Value = '6.518'
$ws.Range("D13").Value = '20.99'
$ws.Range("D15").Value = '7.445'
$ws.Range("D16").Value = '1.001'
$ws.Range("D17").Value = '0.00001129'
$ws.Range("D18").Value = '93.03'
$ws.Range("D20").Value = '17.57'
$ws.Range("D21").Value = '1.000'
$ws.Range("D22").Value = '6.167'
$ws.Range("D25").Value = '2.288'
$ws.Range("D26").Value = '2.547'
$ws.Range("D28").Value = '21.23'
$ws.Range("D29").Value = '158.39'
$ws.Range("D30").Value = '127.67'
$ws.Range("D31").Value = '0.1062'
$ws.Range("D32").Value = '1.069'
$ws.Range("D33").Value = '5.644'
$ws.Range("D34").Value = '3.622'
$ws.Range("D35").Value = '9.547'
$ws.Range("D36").Value = '0.06729'
$ws.Range("D37").Value = '0.02389'
$ws.Range("D38").Value = '0.2193'
$ws.Range("D39").Value = '11.52'
$ws.Range("D40").Value = '0.6383'
$ws.Range("D41").Value = '4.982'
$ws.Range("D42").Value = '1.184'
$ws.Range("D43").Value = '1.0000'
$ws.Range("D44").Value = '13.52'
$ws.Range("D45").Value = '0.6034'
$ws.Range("D46").Value = '3.661'
$ws.Range("D47").Value = '1.259'
$ws.Range("D48").Value = '1.999'
$ws.Range("D49").Value = '124.15'
$ws.Range("D50").Value = '1.198'
$ws.Range("D51").Value = '0.06844'

# Apply the remaining (non-ambiguous) text values: coin names, links, prices, and % changes
$ws.Range("D2").Value = '28.224.80'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '1.870.76'
$ws.Range("E3").Value = '  +3.72%  '
$ws.Range("E4").Value = '  -0.93%  '
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("E7").Value = '  +2.37%  '
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("E9").Value = '  +4.37%  '
$ws.Range("E10").Value = '  +4.35%  '
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("E13").Value = '  +2.12%  '
$ws.Range("D14").Value = '1.877.12'
$ws.Range("E14").Value = '  +3.25%  '
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E20").Value = '  +3.16%  '
$ws.Range("E21").Value = '  -0.91%  '
$ws.Range("E22").Value = '  +3.65%  '
$ws.Range("D23").Value = '28.280.07'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("E24").Value = '  +3.25%  '
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E26").Value = '  +7.82%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.084.73'
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E28").Value = '  +3.80%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E35").Value = '  +7.57%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E38").Value = '  +3.07%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("E39").Value = '  +1.75%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("E40").Value = '  +4.45%  '
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E42").Value = '  +3.25%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E44").Value = '  +3.06%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E47").Value = '  -2.57%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E48").Value = '  +2.94%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E51").Value = '  +1.33%  '
